$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D9 to a formula summing 2225+1466
$ws.Range("D9").Formula = "=2225+1466"

# Update E9 and F9 values
$ws.Range("E9").Value = 2225
$ws.Range("F9").Value = 2225

# Update selection to E10
$ws.Range("E10").Select()
